$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old B->C, old C->D)
$ws.Columns("B").Insert()
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth

# New header for inserted column
$ws.Range("B1").Value = "Strain"

# Fill the Strain values (derived from Subject ID prefix)
$ws.Range("B2").Value = "8034x8043"
$ws.Range("B3").Value = "8034x8043"
$ws.Range("B4").Value = "8034x8043"
$ws.Range("B5").Value = "15119x16521"
$ws.Range("B6").Value = "15119x16521"
$ws.Range("B7").Value = "15119x16521"

$ws.Range("B8").Select()
